$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.77"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'22.94"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'6.352"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.06257"
$ws.Range("G5").Value = "'8"
$ws.Range("D6").Value = "'3.665"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'6.687"
$ws.Range("G7").Value = "'8"
$ws.Range("D8").Value = "'1.356"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'0.8312"
$ws.Range("G9").Value = "'8"
$ws.Range("D10").Value = "'0.01375"
$ws.Range("G10").Value = "'8"
$ws.Range("D11").Value = "'0.1630"
$ws.Range("G11").Value = "'8"
$ws.Range("D12").Value = "'0.08377"
$ws.Range("G12").Value = "'8"
$ws.Range("D13").Value = "'0.03436"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.03111"
$ws.Range("G14").Value = "'8"
$ws.Range("D15").Value = "'0.09312"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'3.875"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'0.001644"
$ws.Range("G17").Value = "'8"
$ws.Range("D18").Value = "'0.04761"
$ws.Range("G18").Value = "'8"
$ws.Range("D19").Value = "'0.006373"
$ws.Range("G19").Value = "'8"
$ws.Range("D20").Value = "'0.005557"
$ws.Range("E20").Value = "'19HotbitTokenHTBWorstin24h"
$ws.Range("G20").Value = "'8"
$ws.Range("D21").Value = "'0.001093"
$ws.Range("G21").Value = "'8"
$ws.Range("G22").Value = "'8"
$ws.Range("D23").Value = "'3.711"
$ws.Range("G23").Value = "'8"
$ws.Range("D24").Value = "'2.322"
$ws.Range("G24").Value = "'8"
$ws.Range("D25").Value = "'0.3342"
$ws.Range("G25").Value = "'8"
$ws.Range("G26").Value = "'8"
$ws.Range("D27").Value = "'0.0002680"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("G38").Value = "'8"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.04704"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.007037"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1165"
$ws.Range("G42").Value = "'8"
$ws.Range("D43").Value = "'0.003350"
$ws.Range("E43").Value = "'42CEJICEJI"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.01216"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.00006286"
$ws.Range("G45").Value = "'8"
$ws.Range("G46").Value = "'8"
$ws.Range("D47").Value = "'0.8999"
$ws.Range("G47").Value = "'8"
$ws.Range("D48").Value = "'0.03438"
$ws.Range("G48").Value = "'8"
$ws.Range("D49").Value = "'0.00002200"
$ws.Range("G49").Value = "'8"
$ws.Range("G50").Value = "'8"
$ws.Range("G51").Value = "'8"
